# Scheduled market-data refresh for the Excalibur leve-profit workbook.
# Pulls fresh Universalis average-price snapshots and rewrites each sheet's
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N) row by row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Horn Glue
$ws.Range("H40").Value = 2061.4375
$ws.Range("J40").Value = 2125.818
$ws.Range("L40").Value = 2125.818
$ws.Range("N40").Value = -2475.818

# Row 97: Potent Spiritbond Potion
$ws.Range("H97").Value = 46374004
$ws.Range("I97").Value = 42929908
$ws.Range("J97").Value = 55558260
$ws.Range("K97").Value = 128789724
$ws.Range("L97").Value = 166674780
$ws.Range("M97").Value = -128789228
$ws.Range("N97").Value = -166675772

# Row 112: Superior Spiritbond Potion
$ws.Range("H112").Value = 6293.3076
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 6485.04
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 19455.12
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -21671.12

# Row 137: Magnesia Whetstone
$ws.Range("H137").Value = 1479218.5
$ws.Range("I137").Value = 2670.6667
$ws.Range("J137").Value = 3251076
$ws.Range("K137").Value = 8012.000100000001
$ws.Range("L137").Value = 9753228
$ws.Range("M137").Value = -5462.000100000001
$ws.Range("N137").Value = -9758328

# Row 138: Cunning Craftsman's Tisane
$ws.Range("H138").Value = 4374.6284
$ws.Range("I138").Value = 2446.111
$ws.Range("J138").Value = 5042.1924
$ws.Range("K138").Value = 7338.333
$ws.Range("L138").Value = 15126.5772
$ws.Range("M138").Value = -2198.333
$ws.Range("N138").Value = -25406.5772

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Steel Ingot
$ws.Range("H32").Value = 2312.5833
$ws.Range("I32").Value = 2217.025
$ws.Range("K32").Value = 2217.025
$ws.Range("M32").Value = -1930.025

# Row 45: Mythril Ingot
$ws.Range("H45").Value = 3153.0833
$ws.Range("I45").Value = 2203
$ws.Range("K45").Value = 2203
$ws.Range("M45").Value = -1826

# Row 132: Mountain Chromite Ingot
$ws.Range("H132").Value = 573292.7
$ws.Range("I132").Value = 716101.6
$ws.Range("J132").Value = 2056.8572
$ws.Range("K132").Value = 2148304.8
$ws.Range("L132").Value = 6170.571599999999
$ws.Range("M132").Value = -2145774.8
$ws.Range("N132").Value = -11230.5716

# Row 138: Titanium Gold Helm of Casting
$ws.Range("H138").Value = 83666
$ws.Range("J138").Value = 83666
$ws.Range("L138").Value = 83666
$ws.Range("N138").Value = -93946

$ws = $wb.Worksheets.Item("BSM")
# Row 80: Titanium Ingot
$ws.Range("H80").Value = 2279.4119
$ws.Range("I80").Value = 920.75
$ws.Range("J80").Value = 2697.4614
$ws.Range("K80").Value = 920.75
$ws.Range("L80").Value = 2697.4614
$ws.Range("M80").Value = 77.25
$ws.Range("N80").Value = -4693.4614

# Row 83: Titanium Ingot
$ws.Range("H83").Value = 2279.4119
$ws.Range("I83").Value = 920.75
$ws.Range("J83").Value = 2697.4614
$ws.Range("K83").Value = 4603.75
$ws.Range("L83").Value = 13487.307
$ws.Range("M83").Value = 388.25
$ws.Range("N83").Value = -23471.307

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Walnut Lumber
$ws.Range("H31").Value = 9736.192999999999
$ws.Range("I31").Value = 4686.1333
$ws.Range("J31").Value = 15347.37
$ws.Range("K31").Value = 4686.1333
$ws.Range("L31").Value = 15347.37
$ws.Range("M31").Value = -4391.1333
$ws.Range("N31").Value = -15937.37

# Row 34: Walnut Lumber
$ws.Range("H34").Value = 9736.192999999999
$ws.Range("I34").Value = 4686.1333
$ws.Range("J34").Value = 15347.37
$ws.Range("K34").Value = 4686.1333
$ws.Range("L34").Value = 15347.37
$ws.Range("M34").Value = -4484.1333
$ws.Range("N34").Value = -15751.37

# Row 58: Mahogany Lumber
$ws.Range("H58").Value = 539300.0600000001
$ws.Range("I58").Value = 1123942.8
$ws.Range("J58").Value = 3377.6667
$ws.Range("K58").Value = 1123942.8
$ws.Range("L58").Value = 3377.6667
$ws.Range("M58").Value = -1123739.8
$ws.Range("N58").Value = -3783.6667

# Row 132: Ginseng Lumber
$ws.Range("H132").Value = 7893.7
$ws.Range("I132").Value = 1938.0834
$ws.Range("J132").Value = 16827.125
$ws.Range("K132").Value = 5814.2502
$ws.Range("L132").Value = 50481.375
$ws.Range("M132").Value = -3284.2502
$ws.Range("N132").Value = -55541.375

# Row 134: Ceiba Lumber
$ws.Range("H134").Value = 9171.883
$ws.Range("I134").Value = 9171.883
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 27515.649
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -24980.649
$ws.Range("N134").ClearContents()

# Row 136: Dark Mahogany Lumber
$ws.Range("H136").Value = 539300.0600000001
$ws.Range("I136").Value = 1123942.8
$ws.Range("J136").Value = 3377.6667
$ws.Range("K136").Value = 3371828.4
$ws.Range("L136").Value = 10133.0001
$ws.Range("M136").Value = -3369278.4
$ws.Range("N136").Value = -15233.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 92: Gyr Abanian Flour
$ws.Range("H92").Value = 600.25
$ws.Range("I92").Value = 380.4
$ws.Range("J92").Value = 966.6667
$ws.Range("K92").Value = 1141.2
$ws.Range("L92").Value = 2900.0001
$ws.Range("M92").Value = 106.8000000000002
$ws.Range("N92").Value = -5396.0001

# Row 132: Cooking Mezcal
$ws.Range("H132").Value = 618.8333
$ws.Range("I132").Value = 565
$ws.Range("J132").Value = 888
$ws.Range("K132").Value = 5085
$ws.Range("L132").Value = 7992
$ws.Range("M132").Value = -2555
$ws.Range("N132").Value = -13052

$ws = $wb.Worksheets.Item("GSM")
# Row 6: Bone Staff
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 5000
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -5226

# Row 16: Decorated Bone Staff
$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5000
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -5500

# Row 132: Lar Ingot
$ws.Range("H132").Value = 378974.56
$ws.Range("I132").Value = 417902.47
$ws.Range("J132").Value = 2671.3333
$ws.Range("K132").Value = 1253707.41
$ws.Range("L132").Value = 8013.999899999999
$ws.Range("M132").Value = -1251177.41
$ws.Range("N132").Value = -13073.9999

$ws = $wb.Worksheets.Item("LTW")
# Row 47: Boarskin Harness
$ws.Range("H47").Value = 26000
$ws.Range("I47").Value = 23000
$ws.Range("J47").Value = 27500
$ws.Range("K47").Value = 23000
$ws.Range("L47").Value = 27500
$ws.Range("M47").Value = -22510
$ws.Range("N47").Value = -28480

# Row 52: Boarskin Harness
$ws.Range("H52").Value = 26000
$ws.Range("I52").Value = 23000
$ws.Range("J52").Value = 27500
$ws.Range("K52").Value = 23000
$ws.Range("L52").Value = 27500
$ws.Range("M52").Value = -22767
$ws.Range("N52").Value = -27966

# Row 122: Gaja Leather
$ws.Range("H122").Value = 4583.3213
$ws.Range("I122").Value = 4254
$ws.Range("K122").Value = 12762
$ws.Range("M122").Value = -10312

# Row 136: Br'aax Leather
$ws.Range("H136").Value = 3044.7058
$ws.Range("I136").Value = 2127.7693
$ws.Range("J136").Value = 6024.75
$ws.Range("K136").Value = 6383.3079
$ws.Range("L136").Value = 18074.25
$ws.Range("M136").Value = -3833.3079
$ws.Range("N136").Value = -23174.25

$ws = $wb.Worksheets.Item("WVR")
# Row 132: Snow Cotton Cloth
$ws.Range("H132").Value = 6495660.5
$ws.Range("I132").Value = 6942602.5
$ws.Range("K132").Value = 20827807.5
$ws.Range("M132").Value = -20825277.5

# Row 136: Sarcenet Cloth
$ws.Range("H136").Value = 5883.147
$ws.Range("I136").Value = 5750.8438
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 17252.5314
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -14702.5314
$ws.Range("N136").Value = -29100
